$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Paragraph "The resutls computed from cvx ..." : demote it out of
#    the numbered list (ilvl 1 / numId 2  ->  ilvl 0 / numId 0), give
#    it a manual left indent of 720 twips with a (zero) hanging
#    indent, and prefix the text with a standalone run containing a
#    single space character.
# ------------------------------------------------------------------
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("The resutls computed from cvx")) {
        $targetPara = $p
        break
    }
}

$targetPara.Range.ListFormat.RemoveNumbers()
$targetPara.LeftIndent = 36        # 36 pt  == 720 twips
$targetPara.FirstLineIndent = -0.0001   # rounds to 0 twips but keeps it a "hanging" indent, not a "firstLine" indent

$pStart = $targetPara.Range.Start
$insPoint = $d.Range($pStart, $pStart)
$insPoint.InsertBefore(" ")

# Nudge formatting on the new space-only run and back again so it does
# not get silently coalesced with the following run when saved.
$spaceRun = $d.Range($pStart, $pStart + 1)
$spaceRun.Font.Bold = 1
$spaceRun.Font.Bold = 0

# ------------------------------------------------------------------
# 2) Paragraph "Debug using KF. Notice two things to investigate: ..."
#    The second sentence is replaced by two new sentences, which stay
#    split across two runs (matching the target XML run layout).
# ------------------------------------------------------------------
$kfPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Debug using KF.")) {
        $kfPara = $p
        break
    }
}

$kfParaStart = $kfPara.Range.Start
$prefix = "Debug using KF. "
$oldTail = "Notice two things to investigate: why in the first step the optimal solution is not the initial solution (in fact, a lot of sensor placement seem to have similar result, as long as the target is in FOV. Check what gradient is used). Also think about a simple case to test the algorithm on."

$sentenceA = "Check the gradient to understand if the infeasibility issue. "
$sentenceB = "Add feasibility recovery function."

$tailRange = $d.Range($kfParaStart + $prefix.Length, $kfParaStart + $prefix.Length + $oldTail.Length)
$origColor = $tailRange.Font.Color
$tailRange.Text = $sentenceA + $sentenceB

# keep this new run visually distinct from "Debug using KF. " for a
# moment, then restore its original colour, so the two stay separate
# runs instead of merging back into the previous one on save.
$newTailRange = $d.Range($kfParaStart + $prefix.Length, $kfParaStart + $prefix.Length + $sentenceA.Length + $sentenceB.Length)
$newTailRange.Font.Color = 1
$newTailRange.Font.Color = $origColor

# Now split sentenceA / sentenceB into two distinct runs the same way.
$splitPos = $kfParaStart + $prefix.Length + $sentenceA.Length
$sentenceBRange = $d.Range($splitPos, $kfParaStart + $prefix.Length + $sentenceA.Length + $sentenceB.Length)
$sentenceBRange.Font.Color = 1
$sentenceBRange.Font.Color = $origColor

# ------------------------------------------------------------------
# 3) "Debug using PF" -> "Debug using PF."
# ------------------------------------------------------------------
$d.Content.Find.Execute("Debug using PF", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Debug using PF.", 2)

Write-Host "edits applied"
